# Updates the cryptos price table (rows 2-51) to the latest snapshot values.
# Mirrors GitHub Actions commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look numeric to Excel (e.g. "1.00", "593.76") but must stay as
# literal text to match the source data (it uses dotted thousands separators elsewhere,
# e.g. "66.741.37"), so force text format before writing those cells.
$forceTextCells = @("D4","D5","D6","D7","D9","D10","D11","D12","D13","D14","D16","D20","D22","D23","D24","D26","D28","D29","D30","D31","D32","D33","D34","D37","D39","D40","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.741.37'
$ws.Range("E2").Value = '  +2.88%  '
$ws.Range("D3").Value = '3.432.03'
$ws.Range("E3").Value = '  +1.62%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '569.13'
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("D6").Value = '183.44'
$ws.Range("E6").Value = '  +5.31%  '
$ws.Range("D7").Value = '0.633'
$ws.Range("E7").Value = '  +1.71%  '
$ws.Range("D8").Value = '3.427.44'
$ws.Range("E8").Value = '  +1.79%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '0.176'
$ws.Range("E10").Value = '  +6.21%  '
$ws.Range("D11").Value = '0.643'
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("D12").Value = '55.08'
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("D13").Value = '0.0000280'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = '9.35'
$ws.Range("E14").Value = '  +2.92%  '
$ws.Range("D15").Value = '3.980.14'
$ws.Range("E15").Value = '  +1.61%  '
$ws.Range("D16").Value = '18.48'
$ws.Range("E16").Value = '  +1.32%  '
$ws.Range("D17").Value = '3.436.60'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").Value = '66.643.01'
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("D20").Value = '12.03'
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("D22").Value = '469.59'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '4.99'
$ws.Range("E23").Value = '  +1.93%  '
$ws.Range("D24").Value = '14.87'
$ws.Range("E24").Value = '  +9.78%  '
$ws.Range("E25").Value = '  +1.50%  '
$ws.Range("D26").Value = '89.46'
$ws.Range("E26").Value = '  +3.25%  '
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("D28").Value = '10.92'
$ws.Range("E28").Value = '  +1.14%  '
$ws.Range("D29").Value = '8.89'
$ws.Range("E29").Value = '  +1.88%  '
$ws.Range("D30").Value = '31.44'
$ws.Range("E30").Value = '  +2.56%  '
$ws.Range("D31").Value = '6.97'
$ws.Range("E31").Value = '  +3.21%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '593.76'
$ws.Range("E32").Value = '  +3.95%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '11.60'
$ws.Range("E33").Value = '  +1.30%  '
$ws.Range("D34").Value = '62.70'
$ws.Range("E34").Value = '  +2.18%  '
$ws.Range("E35").Value = '  +1.81%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = '0.148'
$ws.Range("E37").Value = '  +6.32%  '
$ws.Range("E38").Value = '  +1.02%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = '0.389'
$ws.Range("E39").Value = '  +5.60%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").Value = '36.56'
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("D41").Value = '0.0₃0764'
$ws.Range("E41").Value = '  +3.39%  '
$ws.Range("D42").Value = '3.124.11'
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").Value = '2.91'
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.58'
$ws.Range("E44").Value = '  +5.57%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0425'
$ws.Range("E45").Value = '  +2.66%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").Value = '  +20.76%  '
$ws.Range("D47").Value = '3.24'
$ws.Range("E47").Value = '  +3.39%  '
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '141.66'
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").Value = '8.64'
$ws.Range("E51").Value = '  +4.45%  '
